# The deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml  -> "Office Theme" (clrScheme "Office")   -- currently unused leftover,
#                              only wired to the Notes Master.
#   ppt/theme/theme2.xml  -> "Integral" (clrScheme "Red Violet")   -- the live theme, wired to the
#                              single Slide Master / the presentation's main design.
#
# The target edit swaps the colour palette that is actually applied to the presentation:
# the live design (theme2.xml) should carry the default "Office Theme" colours instead of the
# "Integral" / "Red Violet" ones. We do this the same way a user would from the Design tab --
# by editing the 12 theme colour slots (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) through the
# PowerPoint colour-scheme object model, which writes straight into the clrScheme of the theme
# part that backs the slide master.

function HexToColorRef([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# MsoThemeColorSchemeIndex order: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2,
# 7 accent3, 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink.
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = HexToColorRef $officeThemeColors[$i - 1]
}
